$wb = $excel.ActiveWorkbook

# --- Sheet "person matches": rows 14-18 B/C columns ---
$ws1 = $wb.Worksheets.Item("person matches")
$ws1.Range("C14").Value = "P3214"
$ws1.Range("B15").Value = "eft:munivarman"
$ws1.Range("C15").Value = "P8261"
$ws1.Range("B16").Value = "eft:prajnavarman"
$ws1.Range("C16").Value = "P2548"
$ws1.Range("B17").Value = "eft:dpal-dbyangs"
$ws1.Range("C17").Value = "P8260"
$ws1.Range("B18").Value = "eft:ska-ba-dpal-brtsegs"
$ws1.Range("C18").Value = "P8182"

# --- Sheet "grouped matches": rows 2-55 B/C columns reshuffled, row 56 deleted ---
$ws2 = $wb.Worksheets.Item("grouped matches")
$ws2.Range("B2").Value = "P0TMP092"
$ws2.Range("C2").Value = "{'eft:anandasri-s-'}"
$ws2.Range("B3").Value = "P8263"
$ws2.Range("C3").Value = "{'eft:leki-d-'}"
$ws2.Range("B5").Value = "P8268"
$ws2.Range("C5").Value = "{'eft:buddhaprabha'}"
$ws2.Range("B6").Value = "P8205"
$ws2.Range("C6").Value = "{'eft:band-yesh-de', 'eft:yesh-d-', 'eft:band-yesh-d-', 'eft:zhang-yesh-d-', 'eft:yesh-d-ye-shes-sde-', 'eft:ye-shes-sde'}"
$ws2.Range("B7").Value = "P4CZ16819"
$ws2.Range("C7").Value = "{'eft:sakyaprabha'}"
$ws2.Range("B8").Value = "P8171"
$ws2.Range("C8").Value = "{'eft:dharmasribhadra'}"
$ws2.Range("B9").Value = "P4258"
$ws2.Range("C9").Value = "{'eft:dpal-byor'}"
$ws2.Range("B10").Value = "P8269"
$ws2.Range("C10").Value = "{'eft:dgon-gling-rma'}"
$ws2.Range("B11").Value = "P0TMP080"
$ws2.Range("C11").Value = "{'eft:hwa-shang-zab-mo'}"
$ws2.Range("B12").Value = "P2956"
$ws2.Range("C12").Value = "{'eft:krsnapandita'}"
$ws2.Range("B13").Value = "P8265"
$ws2.Range("C13").Value = "{'eft:ratnaraksita'}"
$ws2.Range("B14").Value = "P8093"
$ws2.Range("C14").Value = "{'eft:kamalagupta'}"
$ws2.Range("B15").Value = "P3456"
$ws2.Range("C15").Value = "{'eft:tshul-khrims-rgyal-ba'}"
$ws2.Range("B16").Value = "P4CZ15137"
$ws2.Range("C16").Value = "{'eft:kumarakalasa'}"
$ws2.Range("B17").Value = "P4255"
$ws2.Range("C17").Value = "{'eft:t-jnanagarbha', 'eft:yesh-nyingpo', 'eft:ye-shes-snying-po'}"
$ws2.Range("B18").Value = "P8266"
$ws2.Range("C18").Value = "{'eft:ch-nyi-tsultrim', 'eft:dharmatasila'}"
$ws2.Range("B19").Value = "P3285"
$ws2.Range("C19").Value = "{'eft:sakya-yesh-'}"
$ws2.Range("B20").Value = "P4242"
$ws2.Range("C20").Value = "{'eft:sherab-lekpa'}"
$ws2.Range("B21").Value = "P4263"
$ws2.Range("C21").Value = "{'eft:dge-ba-dpal'}"
$ws2.Range("B22").Value = "P4259"
$ws2.Range("C22").Value = "{'eft:dpal-gyi-lhun-po', 'eft:ban-de-dpal-gyi-lhun-po', 'eft:palgyi-lh-npo'}"
$ws2.Range("B23").Value = "P2637"
$ws2.Range("C23").Value = "{'eft:trakpa-gyaltsen'}"
$ws2.Range("B24").Value = "P3214"
$ws2.Range("C24").Value = "{'eft:danasila'}"
$ws2.Range("B25").Value = "P5651"
$ws2.Range("C25").Value = "{'eft:pa-tshab-nyi-ma-grags'}"
$ws2.Range("B26").Value = "P8267"
$ws2.Range("C26").Value = "{'eft:vijayasila'}"
$ws2.Range("B27").Value = "P0TMPT007"
$ws2.Range("C27").Value = "{'eft:rnam-par-mi-rtog-pa'}"
$ws2.Range("B28").Value = "P0RK8"
$ws2.Range("C28").Value = "{'eft:dharmapala'}"
$ws2.Range("B29").Value = "P8211"
$ws2.Range("C29").Value = "{'eft:vidyakaraprabha'}"
$ws2.Range("B30").Value = "P8249"
$ws2.Range("C30").Value = "{'eft:dharmakara'}"
$ws2.Range("B31").Value = "P8183"
$ws2.Range("C31").Value = "{'eft:cog-ro-klu-i-rgyal-mtshan', 'eft:klu-i-rgyal-mtshan'}"
$ws2.Range("B32").Value = "P8182"
$ws2.Range("C32").Value = "{'eft:paltsek', 'eft:ban-de-dpal-brtsegs', 'eft:dpal-brtsegs', 'eft:ska-ba-dpal-brtsegs', 'eft:kawa-paltsek-under-the-name-paltsek-raksita-'}"
$ws2.Range("B33").Value = "P8261"
$ws2.Range("C33").Value = "{'eft:munivarman', 'eft:munivarma'}"
$ws2.Range("B34").Value = "P1KG8854"
$ws2.Range("C34").Value = "{'eft:srilendrabodhi', 'eft:surendrabodhi', 'eft:silendrabodhi'}"
$ws2.Range("B35").Value = "P753"
$ws2.Range("C35").Value = "{'eft:rin-chen-bzang-po'}"
$ws2.Range("B37").Value = "P8213"
$ws2.Range("C37").Value = "{'eft:t-vidyakarasimha', 'eft:vidyakarasimha'}"
$ws2.Range("B38").Value = "P00KG07267"
$ws2.Range("C38").Value = "{'eft:sarvajnadeva', 'eft:sarvanyadeva'}"
$ws2.Range("B39").Value = "P3379"
$ws2.Range("C39").Value = "{'eft:dipamkara-srijnana', 'eft:dipamkarasrijnana'}"
$ws2.Range("B40").Value = "P8222"
$ws2.Range("C40").Value = "{'eft:jnanasidhi', 'eft:jnanasiddhi'}"
$ws2.Range("B41").Value = "P0TMP104"
$ws2.Range("C41").Value = "{'eft:punyasambhava'}"
$ws2.Range("B42").Value = "P8219"
$ws2.Range("C42").Value = "{'eft:visuddhasimha'}"
$ws2.Range("B43").Value = "https://lod.dila.edu.tw/resource.php?id=A000089"
$ws2.Range("C43").Value = "{'eft:siladharma'}"
$ws2.Range("B44").Value = "?"
$ws2.Range("C44").Value = "{'eft:sakyasena'}"
$ws2.Range("B45").Value = "P3709"
$ws2.Range("C45").Value = "{'eft:phakpa-sherab'}"
$ws2.Range("B46").Value = "P0TMP098"
$ws2.Range("C46").Value = "{'eft:jinavara'}"
$ws2.Range("B47").Value = "P8273"
$ws2.Range("C47").Value = "{'eft:rin-chen-tsho', 'eft:rinchen-tso'}"
$ws2.Range("B48").Value = "P8151"
$ws2.Range("C48").Value = "{'eft:gayadhara'}"
$ws2.Range("B49").Value = "P8228"
$ws2.Range("C49").Value = "{'eft:surendrabodhi'}"
$ws2.Range("B50").Value = "P2548"
$ws2.Range("C50").Value = "{'eft:prajnavarman', 'eft:prajnavarma'}"
$ws2.Range("B51").Value = "P8217"
$ws2.Range("C51").Value = "{'eft:jnanagarbha', 'eft:t-jnanagarbha'}"
$ws2.Range("B52").Value = "P4CZ16780"
$ws2.Range("C52").Value = "{'eft:manjusrigarbha'}"
$ws2.Range("B53").Value = "P8245"
$ws2.Range("C53").Value = "{'eft:buddhakaravarma'}"
$ws2.Range("B54").Value = "P8260"
$ws2.Range("C54").Value = "{'eft:dpal-dbyangs'}"
$ws2.Range("B55").Value = "P8209"
$ws2.Range("C55").Value = "{'eft:jinamitra', 'eft:dzi-na-mi-tra-k-', 'eft:jinamitra-k-'}"

# Remove the now-empty trailing row (was row 56; content shifted up by reassignment above)
$ws2.Rows.Item(56).Delete()
